# Workbook/worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fossil")

# --- Update numeric values (columns A, D, E, F) ---
# Row 2
$ws.Range("A2").Value = 50
$ws.Range("E2").Value = 409.1
$ws.Range("F2").Value = 412.3

# Row 3
$ws.Range("A3").Value = 21
$ws.Range("D3").Value = 383.5
$ws.Range("E3").Value = 383
$ws.Range("F3").Value = 384

# Row 4
$ws.Range("A4").Value = 36
$ws.Range("D4").Value = 232.15
$ws.Range("E4").Value = 227.3
$ws.Range("F4").Value = 237

# Row 5
$ws.Range("A5").Value = 29
$ws.Range("D5").Value = 232.15
$ws.Range("E5").Value = 227.3
$ws.Range("F5").Value = 237

# Row 6
$ws.Range("A6").Value = 32
$ws.Range("D6").Value = 232.15
$ws.Range("E6").Value = 227.3
$ws.Range("F6").Value = 237

# Row 7
$ws.Range("A7").Value = 79
$ws.Range("D7").Value = 383.5
$ws.Range("E7").Value = 383.5

# Row 8
$ws.Range("A8").Value = 53

# --- Apply formatting ---
# Bold font on A2:A3
$ws.Range("A2").Font.Bold = $true
$ws.Range("A3").Font.Bold = $true

# Yellow fill on D7:D8
$ws.Range("D7").Interior.Color = 65535
$ws.Range("D8").Interior.Color = 65535

# --- Update selection to D7 ---
$ws.Range("D7").Select() | Out-Null
